$d = $word.ActiveDocument

function Split-Into-Runs {
    param(
        [string]$OldText,
        [string[]]$Parts
    )

    $range = $d.Content
    $found = $range.Find.Execute($OldText, $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $OldText"
    }

    $start = $range.Start

    # Force a run boundary at each word/space boundary by adding a
    # temporary bookmark around each part, then deleting the bookmarks.
    # This splits the single run into one run per part without leaving
    # any stray run-formatting (rPr) behind.
    $pos = $start
    $bookmarkNames = @()
    $i = 0
    foreach ($part in $Parts) {
        $len = $part.Length
        $i = $i + 1
        $bmName = "zzTmpSplit$i"
        $r = $d.Range($pos, $pos + $len)
        $d.Bookmarks.Add($bmName, $r) | Out-Null
        $bookmarkNames += $bmName
        $pos = $pos + $len
    }

    foreach ($bmName in $bookmarkNames) {
        $d.Bookmarks($bmName).Delete()
    }
}

Split-Into-Runs "Questions: Trigonometry (degrees)" @("Questions:", " ", "Trigonometry", " ", "(degrees)")
Split-Into-Runs "A selection of questions on trigonometry, where angles are measured in degrees." @("A", " ", "selection", " ", "of", " ", "questions", " ", "on", " ", "trigonometry,", " ", "where", " ", "angles", " ", "are", " ", "measured", " ", "in", " ", "degrees.")
